$d = $word.ActiveDocument
$d.Content.Find.Execute("Cybersecurity, Networking", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cybersecurity, Networking", 2)
